$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 423, shifting existing rows 423:500 down to 424:501.
$ws.Rows.Item(423).Insert()

# Populate the newly inserted row 423 with the new data record.
$ws.Range("A423").Value = 4
$ws.Range("B423").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C423").Value = "Los Lagos"
$ws.Range("D423").Value = 45258
$ws.Range("E423").Value = 10
$ws.Range("F423").Value = "Fruta"
$ws.Range("G423").Value = 100108
$ws.Range("H423").Value = "Tropicales y subtropicales"
$ws.Range("I423").Value = 100108005
$ws.Range("J423").Value = "Piña"
$ws.Range("K423").Value = "Caramelo"
$ws.Range("L423").Value = "Primera"
$ws.Range("M423").Value = 200
$ws.Range("N423").Value = 26000
$ws.Range("O423").Value = 26000
$ws.Range("P423").Value = 26000
$ws.Range("Q423").Value = "`$/caja 12 unidades"
$ws.Range("R423").Value = "Ecuador"
$ws.Range("S423").Value = 2167
$ws.Range("T423").Value = 12
